# Workbook: municipal area table for "მარტვილი" (Martvili) municipality.
# The sheet is reshaped from a 3-year (1989/2002/2014) comparison table into
# a single-year (2014) snapshot, the population-census caption row is
# dropped, and the worksheet/tab is given its proper Georgian name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from the generic "1" to "მარტვილი".
$ws.Name = "მარტვილი"

# Drop the 2002/1989 columns (C:D) - only the 2014 figures are kept.
$ws.Range("C1:D1").EntireColumn.Delete()

# Drop the "(მოსახლეობის აღწერის შედეგებით)" caption row (old row 2).
$ws.Range("A2").EntireRow.Delete()

# --- Rewrite remaining content to match the simplified layout ---
# Row 1: title
$ws.Range("A1").Value2 = "მარტვილის მუნიციპალიტეტის ფართობი"
$ws.Range("B1").Clear()

# Row 2: now a blank spacer row
$ws.Range("B2").Clear()

# Row 3: unit caption
$ws.Range("A3").Value2 = "(კვ. კმ)"

# Row 4: year header (now just 2014)
$ws.Range("B4").Value2 = 2014

# Row 5: area label + value
$ws.Range("A5").Value2 = "ფართობი"
$ws.Range("B5").Value2 = 880.6

# Match the selection left behind in the saved file.
$ws.Range("A2").Select()
